$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# First, copy the cell formatting for every affected row from a matching
# "template" row so that the new cells pick up the correct style indices
# (the sheet alternates between two visual styles for legibility).
#   rows 52 and 54 -> "alternate" style (same as existing row 50)
#   rows 53, 55, 56, 57 -> "regular" style (same as existing row 51)
# ---------------------------------------------------------------------------
$ws.Range("A50:E50").Copy()
$ws.Range("A52:E52").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A54:E54").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A51:E51").Copy()
$ws.Range("A53:E53").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A55:E55").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A56:E56").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 57 has no file name in column A at all (unlike the other new rows),
# so only copy formatting for columns B:E.
$ws.Range("B51:E51").Copy()
$ws.Range("B57:E57").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Populate the new text cells. The order in which new strings are first
# written determines their position in the shared-string table, so the
# assignments below are deliberately sequenced to reproduce that order.
# ---------------------------------------------------------------------------
$ws.Range("A52").Value = "SCRIPT/T01P02A/um2508.ssb"

$ws.Range("C53").Value = " Oh? You want to know about the\ngrand master of all things bad?"
$ws.Range("C54").Value = " Let me put it as simply...[K]\nEeeeeeek!"

$ws.Range("A53").Value = "SCRIPT/P01P01A/us0101.ssb"

$ws.Range("D53").Value = " Ой? Вы хотите знать о Гранд\nмастере всего самого плохого?"
$ws.Range("D54").Value = " Я просто скажу...[K]\nИиииииии!"

$ws.Range("E53").Value = " Ïê? Âú öïóéóå èîàóû ï Ãñàîä\níàòóåñå âòåãï òàíïãï ðìïöïãï?"
$ws.Range("E54").Value = " Ÿ ðñïòóï òëàçô...[K]\nÉééééééé!"

$ws.Range("A54").Value = "SCRIPT/P01P01A/us0103.ssb"

$ws.Range("C55").Value = " Time passes so quickly!"
$ws.Range("C56").Value = " It seems as if it were only\nyesterday that you two showed up at the guild."
$ws.Range("C57").Value = " The two of you then...[K] Oh, you\nwere adorable! Oh my gosh!"

$ws.Range("A55").Value = "SCRIPT/G01P03A/us0109.ssb"

$ws.Range("D55").Value = " Как же летит время!"
$ws.Range("D56").Value = " Кажется, что вы двое появились\nв гильдии только вчера."
$ws.Range("D57").Value = " Тогда вы были...[K] О, вы были\nтакими очаровашками! О боже мой!"

$ws.Range("E55").Value = " Ëàë çå ìåóéó âñåíÿ!"
$ws.Range("E56").Value = " Ëàçåóòÿ, œóï âú äâïå ðïÿâéìéòû\nâ ãéìûäéé óïìûëï âœåñà."
$ws.Range("E57").Value = " Óïãäà âú áúìé...[K] Ï, âú áúìé\nóàëéíé ïœàñïâàšëàíé! Ï áïçå íïê!"

$ws.Range("A56").Value = "SCRIPT/G01P03A/us3101.ssb "

# ---------------------------------------------------------------------------
# Numeric "line number" cells (column B). These are plain numbers, not
# shared strings, so their order has no effect on the string table.
# ---------------------------------------------------------------------------
$ws.Range("B52").Value = 376
$ws.Range("B53").Value = 351
$ws.Range("B54").Value = 354
$ws.Range("B55").Value = 307
$ws.Range("B56").Value = 310
$ws.Range("B57").Value = 313

# ---------------------------------------------------------------------------
# Row heights.
# ---------------------------------------------------------------------------
$ws.Rows.Item(52).RowHeight = 43.2
$ws.Rows.Item(53).RowHeight = 32.4
$ws.Rows.Item(54).RowHeight = 30
$ws.Rows.Item(55).RowHeight = 43.2
$ws.Rows.Item(56).RowHeight = 43.2
$ws.Rows.Item(57).RowHeight = 21.6

# ---------------------------------------------------------------------------
# Update sheet view to focus near the newly-added rows.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 53
$win.ScrollColumn = 1
$ws.Range("C55").Select()
